# "Generate Report for Handback"
# The 6530e492... file was previously "Ready for handoff". It has now been
# handed back and is in sync with en-US, same as the c9e0a3ca... file.
# This updates the per-language status rows (and Overview roll-up) and
# records the new handback timestamps, keeping each language sheet's
# rows ordered with the most-recently-handled-back file first.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: both files are now "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

foreach ($h in $overview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "6530e492-db13-4656-95d7-bf6b0e1b5df7.md"
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = "c9e0a3ca-4e8b-4613-9f4c-765827c62fce.md"
    }
}

# ---------------------------------------------------------------------
# Helper data: per-language sheet name -> new "Latest Handback DateTime"
# for the 6530e492... file, now that it has been handed back.
# ---------------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; HandbackDateTime = "2016-03-08 14:39:20" },
    @{ Sheet = "de-de"; HandbackDateTime = "2016-03-08 14:39:33" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Row 2 held the c9e0a3ca... file, row 3 held the 6530e492... file.
    # The 6530e492... file is now handed back, so it moves up to row 2
    # (swap the two rows across all used columns A-H).
    foreach ($col in @("A", "B", "C", "D", "E", "F", "G", "H")) {
        $row2Value = $ws.Range($col + "2").Value()
        $row3Value = $ws.Range($col + "3").Value()
        $ws.Range($col + "2").Value = $row3Value
        $ws.Range($col + "3").Value = $row2Value
    }

    # Update status + handback datetime for the newly-promoted row 2
    # (6530e492...) to reflect the handback that just happened.
    $ws.Range("B2").Value = $handedBack
    $ws.Range("G2").Value = $lang.HandbackDateTime

    # Keep hyperlink display text in sync with the swapped cell contents.
    $file6530 = "6530e492-db13-4656-95d7-bf6b0e1b5df7.md"
    $file6530Handback = "6530e492-db13-4656-95d7-bf6b0e1b5df7.40a1b2fde74259dd5f2ea8619ac2fe4bde1bad9d." + $lang.Sheet + ".xlf"
    $fileC9e0 = "c9e0a3ca-4e8b-4613-9f4c-765827c62fce.md"
    $fileC9e0Handback = "c9e0a3ca-4e8b-4613-9f4c-765827c62fce.ae2fec7d088533bfad3fc99329a870ec66f54517." + $lang.Sheet + ".xlf"

    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq '$A$2' -or $addr -eq '$E$2') {
            $h.TextToDisplay = $file6530
        } elseif ($addr -eq '$C$2' -or $addr -eq '$F$2') {
            $h.TextToDisplay = $file6530Handback
        } elseif ($addr -eq '$A$3' -or $addr -eq '$E$3') {
            $h.TextToDisplay = $fileC9e0
        } elseif ($addr -eq '$C$3' -or $addr -eq '$F$3') {
            $h.TextToDisplay = $fileC9e0Handback
        }
    }
}
